# table1.xlsx -- "data and name updates"
#
# The sheet holds excess-mortality summary stats for a handful of
# countries/years (Year, Country, Expected_Mortality, Deaths_num,
# Cum_excess_death, Percent_excess_death in columns A-F).
#
# This refreshes the 2020 rows (Spain/Sweden/Switzerland, rows 10-12)
# with corrected Deaths_num / Cum_excess_death / Percent_excess_death
# figures (and a corrected Expected_Mortality for Switzerland).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - Spain, 2020
$ws.Range("D10").Value = 479760
$ws.Range("E10").Value = 54276
$ws.Range("F10").Value = 12.8

# Row 11 - Sweden, 2020
$ws.Range("D11").Value = 97870
$ws.Range("E11").Value = 8895
$ws.Range("F11").Value = 10

# Row 12 - Switzerland, 2020
$ws.Range("C12").Value = 67445
$ws.Range("D12").Value = 75570
$ws.Range("E12").Value = 8125
$ws.Range("F12").Value = 12
